# Add a new forecast vintage column (BB) to Sheet1.
#
# Column BB mirrors column BA for the historical rows (3-18), while the
# most recent rows (19-21) receive newly computed forecast values. Row 1
# gets a new date header value (BB1), formatted the same way as the rest
# of row 1's header cells (same style as BA1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header date in BB1 ---------------------------------------
$ws.Range("BB1").Value = 45986

# Copy the number format/style from BA1 onto BB1 (paste formats only, so
# the underlying value we just set is preserved).
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)

# --- Rows 3-18: BB repeats the BA value for that row ----------------------
foreach ($r in 3..18) {
    $baCell = $ws.Cells.Item($r, 53)   # column BA = 53
    $bbCell = $ws.Cells.Item($r, 54)   # column BB = 54
    $bbCell.Value = $baCell.Value()
}

# --- Rows 19-21: BB gets newly computed forecast values --------------------
$ws.Range("BB19").Value = -0.7200474048664085
$ws.Range("BB20").Value = -2.181280391105744
$ws.Range("BB21").Value = -2.104371875253941
